{"js": "// \"added part of Oregon\" \u2014 move the \"Cascade Lakes National Scenic Byway\"\n// paragraph earlier in the Bend-area list (right after the \"Bend is\n// world-famous...\" paragraph, before \"Painted Hills & Smith Rock State\n// Park\"), and insert a brand-new short \"Sisters\" bullet right after the\n// \"Misery Ridge Trail\" paragraph (before \"The town of Sisters...\" bullet).\n//\n// (The diff's many <w:proofErr> removals / run-merges scattered across the\n// document are Word's own proofing-pass churn \u2014 they don't change visible\n// text, so we don't need to replicate them programmatically.)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfunction findIndex(predicate) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (predicate(paragraphs.items[i].text)) {\n      return i;\n    }\n  }\n  return -1;\n}\n\nconst cascadeText =\n  \"To the south, visit Cascade Lakes National Scenic Byway. This 66-mile \" +\n  \"drive will take between 3 to 5 hours and is scattered with alpine lakes \" +\n  \"and snowcapped mountains. A favorite hike here is the Ray Atkson Loop \" +\n  \"Trail that is 2.3 miles long and located near Sparks Lake. It is also a \" +\n  \"fun place to bring a kayak and paddle for a while.\";\n\n// Grab a handle on the *original* \"Cascade Lakes\" paragraph before we insert\n// anything, so a later lookup-by-text can't accidentally grab the new copy.\nconst idxOldCascade = findIndex(\n  (t) => t.indexOf(\"To the south, visit Cascade Lakes\") === 0\n);\nif (idxOldCascade === -1) {\n  throw new Error(\"Could not find the original 'Cascade Lakes' paragraph\");\n}\nconst oldCascadeParagraph = paragraphs.items[idxOldCascade];\n\n// 1) Insert the \"Cascade Lakes\" paragraph right after \"Bend is world-famous...\".\nconst idxBend = findIndex((t) => t.indexOf(\"Bend is world-famous\") === 0);\nif (idxBend === -1) {\n  throw new Error(\"Could not find the 'Bend is world-famous...' paragraph\");\n}\nparagraphs.items[idxBend].insertParagraph(cascadeText, Word.InsertLocation.after);\nawait context.sync();\n\n// 2) Insert a brand-new \"Sisters\" bullet right after the \"Misery Ridge Trail\"\n//    paragraph, at the top list level (ilvl 0), matching \"Bend\" / \"Painted\n//    Hills & Smith Rock State Park\" siblings.\nparagraphs.load(\"text\");\nawait context.sync();\nconst idxMisery = findIndex((t) => t.indexOf(\"Misery Ridge\") !== -1);\nif (idxMisery === -1) {\n  throw new Error(\"Could not find the 'Misery Ridge Trail' paragraph\");\n}\nconst sistersPara = paragraphs.items[idxMisery].insertParagraph(\n  \"Sisters\",\n  Word.InsertLocation.after\n);\nconst sistersListItem = sistersPara.listItemOrNullObject;\nsistersListItem.load(\"level\");\nawait context.sync();\nsistersListItem.level = 0;\nawait context.sync();\n\n// 3) Remove the old \"Cascade Lakes\" paragraph from its original location\n//    (between \"The town of Sisters...\" and \"Rogue Valley Area\"), using the\n//    handle captured before any insertions happened.\noldCascadeParagraph.delete();\nawait context.sync();\n", "ps1": "# \"added part of Oregon\" \u2014 move the \"Cascade Lakes National Scenic Byway\"\n# paragraph earlier in the Bend-area list (right after the \"Bend is\n# world-famous...\" paragraph, before \"Painted Hills & Smith Rock State\n# Park\"), and insert a brand-new short \"Sisters\" bullet right after the\n# \"Misery Ridge Trail\" paragraph (before \"The town of Sisters...\" bullet).\n#\n# (The diff's many proofErr-tag removals / run-merges scattered across the\n# document are Word's own proofing-pass churn \u2014 they don't change visible\n# text, so we don't need to replicate them programmatically.)\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($needle, $startsWith) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        if ($startsWith) {\n            if ($t.StartsWith($needle)) { return $i }\n        } else {\n            if ($t.Contains($needle)) { return $i }\n        }\n    }\n    return -1\n}\n\n$cascadeText = \"To the south, visit Cascade Lakes National Scenic Byway. This 66-mile drive will take between 3 to 5 hours and is scattered with alpine lakes and snowcapped mountains. A favorite hike here is the Ray Atkson Loop Trail that is 2.3 miles long and located near Sparks Lake. It is also a fun place to bring a kayak and paddle for a while.\"\n\n# 1) Remove the \"Cascade Lakes\" paragraph from its original location FIRST\n#    (between \"The town of Sisters...\" and \"Rogue Valley Area\") \u2014 doing the\n#    delete before any inserts avoids any ambiguity about which copy a later\n#    text search would find.\n$idxOldCascade = Find-ParagraphIndex \"To the south, visit Cascade Lakes\" $true\nif ($idxOldCascade -eq -1) {\n    throw \"Could not find the original 'Cascade Lakes' paragraph\"\n}\n$d.Paragraphs.Item($idxOldCascade).Range.Delete()\n\n# 2) Insert the \"Cascade Lakes\" paragraph right after \"Bend is world-famous...\".\n$idxBend = Find-ParagraphIndex \"Bend is world-famous\" $true\nif ($idxBend -eq -1) {\n    throw \"Could not find the 'Bend is world-famous...' paragraph\"\n}\n$d.Paragraphs.Item($idxBend).Range.InsertParagraphAfter()\n$cascadePara = $d.Paragraphs.Item($idxBend + 1)\n$cascadePara.Range.Text = $cascadeText\n\n# 3) Insert a brand-new \"Sisters\" bullet right after the \"Misery Ridge Trail\"\n#    paragraph, at the top list level (1 in COM's 1-based ListLevelNumber,\n#    i.e. ilvl 0), matching \"Bend\" / \"Painted Hills & Smith Rock State Park\"\n#    siblings.\n$idxMisery = Find-ParagraphIndex \"Misery Ridge\" $false\nif ($idxMisery -eq -1) {\n    throw \"Could not find the 'Misery Ridge Trail' paragraph\"\n}\n$d.Paragraphs.Item($idxMisery).Range.InsertParagraphAfter()\n$sistersPara = $d.Paragraphs.Item($idxMisery + 1)\n$sistersPara.Range.Text = \"Sisters\"\n$sistersPara.Range.ListFormat.ListLevelNumber = 1\n"}
